$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" footer field text (the
#    slide master + every slide layout show the last-saved date) from
#    23/01/2019 to 18/02/2019.
# ---------------------------------------------------------------------
function Set-DatePlaceholderText($shapes, $text) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "18/02/2019"
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $cl = $master.CustomLayouts.Item($i)
    Set-DatePlaceholderText $cl.Shapes "18/02/2019"
}

# ---------------------------------------------------------------------
# 2. Append a new blank slide (slide 5) with a reminder note textbox.
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)   # ppLayoutBlank = 12

# Burn through two throw-away shape ids/names so the real textbox lands
# on id=4 / "TextBox 3", matching the author's original editing session.
$junk1 = $newSlide.Shapes.AddTextbox(1, 0, 0, 10, 10)
$junk1.Delete()
$junk2 = $newSlide.Shapes.AddTextbox(1, 0, 0, 10, 10)
$junk2.Delete()

# Shape position/size in EMU -> points (1 pt = 12700 EMU) so the saved
# OOXML off/ext values come out exactly as authored.
$left = 5022574 / 12700
$top = 1696278 / 12700
$width = 3968779 / 12700
$height = 369332 / 12700

$tb = $newSlide.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Fill.Visible = 0
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1

$tr = $tb.TextFrame.TextRange
$tr.Text = "Email Dom the specs for "
$tr.LanguageID = "en-GB"
$tr2 = $tr.InsertAfter("simulation stuff")
$tr2.LanguageID = "en-GB"
